$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the query text in B2 (remove the trailing "Cohort" column from the RETURN clause)
$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`r`n`r`nMATCH (c)<--(diag:diagnosis)`r`n MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`r`n`t`tWHERE s.clinical_study_designation IN ['UBC01'] and demo.sex in ['Female'] and demo.neutered_indicator in ['No']  OPTIONAL MATCH (samp:sample)-->(c)`r`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`r`nWITH DISTINCT c, s, demo, diag, co`r`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`r`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`r`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`r`n        coalesce(demo.breed, '') AS Breed ,`r`n        coalesce(diag.disease_term, '') AS Diagnosis ,`r`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`r`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`r`n        coalesce(demo.sex, '') AS Sex ,`r`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`r`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`r`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

# Update the sheet view: scroll back to column A / row 2, and select B2 instead of C2:C2:F2
$excel.Goto($ws.Range("A2"), $true)
$ws.Range("B2").Select()

# Adjust row 2 height to match the shortened cell content
$ws.Rows(2).RowHeight = 290
